$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (C) column is bumped by one day (46073 -> 46074) for every
# data row, and the order of the case rows (4-11) has been reshuffled as a
# new data row (A 23798-2024) was discovered, shifting several rows. We
# rewrite the Beteckning (A), Datum (B), Förändrad (C) and Area (G) values
# for every data row explicitly to match the refreshed source data.

$rows = @(
    @{ Row = 2;  A = "A 25353-2022"; B = 44732;              C = 46074; G = 1.5 },
    @{ Row = 3;  A = "A 25351-2022"; B = 44732;              C = 46074; G = 2.8 },
    @{ Row = 4;  A = "A 26262-2024"; B = 45468.66077546297;  C = 46074; G = 0.6 },
    @{ Row = 5;  A = "A 14517-2023"; B = 45012;              C = 46074; G = 0.6 },
    @{ Row = 6;  A = "A 50762-2025"; B = 45946;              C = 46074; G = 2.7 },
    @{ Row = 7;  A = "A 14516-2023"; B = 45012.86600694444;  C = 46074; G = 0.4 },
    @{ Row = 8;  A = "A 23798-2024"; B = 45455.43208333333;  C = 46074; G = 1.3 },
    @{ Row = 9;  A = "A 4156-2023";  B = 44953;              C = 46074; G = 1.5 },
    @{ Row = 10; A = "A 8679-2026";  B = 46066;              C = 46074; G = 2.1 },
    @{ Row = 11; A = "A 8929-2026";  B = 46069.34543981482;  C = 46074; G = 1.2 },
    @{ Row = 12; A = "A 4159-2023";  B = 44953;              C = 46074; G = 0.5 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
